# Updated cryptos list — refresh Price (column D) and Volume(1h) (column E)
# values for the coin ranking table on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "28.593.39";  E = "  +1.89%  " },
    @{ Row = 3;  D = "1.580.07";   E = "  +0.34%  " },
    @{ Row = 4;  D = $null;        E = "  +0.31%  " },
    @{ Row = 5;  D = "212.54";     E = "  -0.09%  " },
    @{ Row = 6;  D = "0.492";      E = "  +0.08%  " },
    @{ Row = 7;  D = $null;        E = "  +0.26%  " },
    @{ Row = 8;  D = "46.89";      E = "  +7.82%  " },
    @{ Row = 9;  D = "24.21";      E = "  +4.46%  " },
    @{ Row = 10; D = $null;        E = "  -1.26%  " },
    @{ Row = 11; D = $null;        E = "  -1.02%  " },
    @{ Row = 12; D = $null;        E = "  +0.08%  " },
    @{ Row = 13; D = "1.805.78";   E = "  +0.36%  " },
    @{ Row = 14; D = "1.569.65";   E = "  -0.28%  " },
    @{ Row = 15; D = "0.523";      E = "  +0.31%  " },
    @{ Row = 16; D = $null;        E = "  -1.37%  " },
    @{ Row = 17; D = "28.579.78";  E = "  +1.94%  " },
    @{ Row = 18; D = "62.34";      E = "  -1.85%  " },
    @{ Row = 19; D = "228.99";     E = $null },
    @{ Row = 20; D = "7.40";       E = "  -0.63%  " },
    @{ Row = 21; D = "0.0₃0695";   E = "  -1.67%  " },
    @{ Row = 22; D = $null;        E = "  +0.31%  " },
    @{ Row = 23; D = $null;        E = "  -4.78%  " },
    @{ Row = 24; D = $null;        E = "  -1.80%  " },
    @{ Row = 25; D = $null;        E = "  +5.30%  " },
    @{ Row = 26; D = "151.22";     E = "  -0.80%  " },
    @{ Row = 27; D = "15.00";      E = "  -1.49%  " },
    @{ Row = 28; D = $null;        E = "  -1.58%  " },
    @{ Row = 29; D = $null;        E = "  -1.74%  " },
    @{ Row = 30; D = $null;        E = "  +0.31%  " },
    @{ Row = 31; D = $null;        E = "  -2.03%  " },
    @{ Row = 32; D = "0.0464";     E = "  -1.82%  " },
    @{ Row = 33; D = $null;        E = "  +0.00%  " },
    @{ Row = 34; D = "3.15";       E = "  +0.53%  " },
    @{ Row = 35; D = "1.397.94";   E = "  -1.37%  " },
    @{ Row = 36; D = $null;        E = "  -2.24%  " },
    @{ Row = 37; D = $null;        E = "  -2.65%  " },
    @{ Row = 38; D = $null;        E = "  +1.88%  " },
    @{ Row = 39; D = "2.61";       E = "  +5.12%  " },
    @{ Row = 40; D = $null;        E = "  -0.42%  " },
    @{ Row = 41; D = $null;        E = "  -1.59%  " },
    @{ Row = 42; D = $null;        E = "  +0.35%  " },
    @{ Row = 44; D = "5.62";       E = "  -0.81%  " },
    @{ Row = 45; D = $null;        E = "  +2.63%  " },
    @{ Row = 46; D = $null;        E = "  +0.41%  " },
    @{ Row = 47; D = "62.97";      E = "  -1.25%  " },
    @{ Row = 48; D = "1.717.53";   E = "  +0.29%  " },
    @{ Row = 49; D = "86.03";      E = "  -1.04%  " },
    @{ Row = 50; D = $null;        E = "  -1.91%  " },
    @{ Row = 51; D = $null;        E = "  -1.28%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        # Force the cell to be treated as text so values such as "212.54"
        # or "15.00" are not silently re-interpreted as numbers (which would
        # also strip significant trailing zeros). Restore the original
        # (style-less) appearance immediately afterwards.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
